$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 1450
$ws.Range("I20").Value = 1450
$ws.Range("K20").Value = 1450
$ws.Range("M20").Value = -1220

# Row 35
$ws.Range("H35").Value = 1450
$ws.Range("I35").Value = 1450
$ws.Range("K35").Value = 1450
$ws.Range("M35").Value = -1071

# Row 125
$ws.Range("H125").Value = 2110
$ws.Range("J125").Value = 2688
$ws.Range("L125").Value = 24192
$ws.Range("N125").Value = -29112

# Row 138
$ws.Range("H138").Value = 6329
$ws.Range("J138").Value = 9128.370000000001
$ws.Range("L138").Value = 27385.11
$ws.Range("N138").Value = -37665.11

# Row 141
$ws.Range("H141").Value = 3039.2
$ws.Range("I141").Value = 3349
$ws.Range("J141").Value = 1800
$ws.Range("K141").Value = 10047
$ws.Range("L141").Value = 5400
$ws.Range("M141").Value = -4867
$ws.Range("N141").Value = -15760

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3030.7
$ws.Range("I32").Value = 2923.037
$ws.Range("K32").Value = 2923.037
$ws.Range("M32").Value = -2636.037

# Row 45
$ws.Range("H45").Value = 2449.5
$ws.Range("I45").Value = 2449.5
$ws.Range("K45").Value = 2449.5
$ws.Range("M45").Value = -2072.5

# Row 61
$ws.Range("H61").Value = 9172.5
$ws.Range("I61").Value = 11231
$ws.Range("K61").Value = 11231
$ws.Range("M61").Value = -11019

# Row 63
$ws.Range("H63").Value = 3999.5
$ws.Range("I63").Value = 1999
$ws.Range("K63").Value = 1999
$ws.Range("M63").Value = -1313

# Row 66
$ws.Range("H66").Value = 3999.5
$ws.Range("I66").Value = 1999
$ws.Range("K66").Value = 9995
$ws.Range("M66").Value = -6563

# Row 80
$ws.Range("H80").Value = 78104.5
$ws.Range("J80").Value = 78104.5
$ws.Range("L80").Value = 78104.5
$ws.Range("N80").Value = -80100.5

# Row 83
$ws.Range("H83").Value = 78104.5
$ws.Range("J83").Value = 78104.5
$ws.Range("L83").Value = 234313.5
$ws.Range("N83").Value = -244297.5

# Row 88
$ws.Range("H88").Value = 2999
$ws.Range("I88").Value = 2999
$ws.Range("K88").Value = 2999
$ws.Range("M88").Value = -2593

# Row 91
$ws.Range("H91").Value = 2999
$ws.Range("I91").Value = 2999
$ws.Range("K91").Value = 2999
$ws.Range("M91").Value = -1595

# Row 97
$ws.Range("H97").Value = 1042.9333
$ws.Range("I97").Value = 245.54546
$ws.Range("J97").Value = 3235.75
$ws.Range("K97").Value = 245.54546
$ws.Range("L97").Value = 3235.75
$ws.Range("M97").Value = 250.45454
$ws.Range("N97").Value = -4227.75

# Row 136
$ws.Range("H136").Value = 9172.5
$ws.Range("I136").Value = 11231
$ws.Range("K136").Value = 33693
$ws.Range("M136").Value = -31143

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 45415.668
$ws.Range("I82").Value = 10257
$ws.Range("K82").Value = 10257
$ws.Range("M82").Value = -9874

# Row 85
$ws.Range("H85").Value = 45415.668
$ws.Range("I85").Value = 10257
$ws.Range("K85").Value = 10257
$ws.Range("M85").Value = -8931

# Row 86
$ws.Range("H86").Value = 3936.35
$ws.Range("I86").Value = 3822
$ws.Range("J86").Value = 4279.4
$ws.Range("K86").Value = 3822
$ws.Range("L86").Value = 4279.4
$ws.Range("M86").Value = -2699
$ws.Range("N86").Value = -6525.4

# Row 89
$ws.Range("H89").Value = 3936.35
$ws.Range("I89").Value = 3822
$ws.Range("J89").Value = 4279.4
$ws.Range("K89").Value = 19110
$ws.Range("L89").Value = 21397
$ws.Range("M89").Value = -13494
$ws.Range("N89").Value = -32629

# Row 94
$ws.Range("H94").Value = 500
$ws.Range("I94").Value = 500
$ws.Range("K94").Value = 500
$ws.Range("M94").Value = -49

# Row 105
$ws.Range("H105").Value = 2818.8
$ws.Range("I105").Value = 2818.8
$ws.Range("K105").Value = 2818.8
$ws.Range("M105").Value = -1071.8

$ws = $wb.Worksheets.Item("CRP")
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

# Row 58
$ws.Range("H58").Value = 1782.6154
$ws.Range("I58").Value = 1782.6154
$ws.Range("K58").Value = 1782.6154
$ws.Range("M58").Value = -1579.6154

# Row 68
$ws.Range("H68").Value = 34987.5
$ws.Range("J68").Value = 34987.5
$ws.Range("L68").Value = 34987.5
$ws.Range("N68").Value = -36485.5

# Row 71
$ws.Range("H71").Value = 34987.5
$ws.Range("J71").Value = 34987.5
$ws.Range("L71").Value = 104962.5
$ws.Range("N71").Value = -112450.5

# Row 74
$ws.Range("H74").Value = 34542.25
$ws.Range("J74").Value = 34542.25
$ws.Range("L74").Value = 34542.25
$ws.Range("N74").Value = -36290.25

# Row 77
$ws.Range("H77").Value = 34542.25
$ws.Range("J77").Value = 34542.25
$ws.Range("L77").Value = 103626.75
$ws.Range("N77").Value = -112362.75

# Row 132
$ws.Range("H132").Value = 3851.423
$ws.Range("I132").Value = 3723.4348
$ws.Range("K132").Value = 11170.3044
$ws.Range("M132").Value = -8640.304400000001

# Row 134
$ws.Range("H134").Value = 3134.8462
$ws.Range("I134").Value = 3167.4092
$ws.Range("K134").Value = 9502.2276
$ws.Range("M134").Value = -6967.2276

# Row 136
$ws.Range("H136").Value = 1782.6154
$ws.Range("I136").Value = 1782.6154
$ws.Range("K136").Value = 5347.8462
$ws.Range("M136").Value = -2797.8462

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 335.8889
$ws.Range("I12").Value = 257.8
$ws.Range("K12").Value = 773.4000000000001
$ws.Range("M12").Value = -600.4000000000001

# Row 92
$ws.Range("H92").Value = 996.6667
$ws.Range("I92").Value = 996.6667
$ws.Range("K92").Value = 2990.0001
$ws.Range("M92").Value = -1742.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1376.75
$ws.Range("I97").Value = 1669.3334
$ws.Range("K97").Value = 1669.3334
$ws.Range("M97").Value = -1173.3334

# Row 102
$ws.Range("H102").Value = 3699
$ws.Range("I102").Value = 3699
$ws.Range("K102").Value = 3699
$ws.Range("M102").Value = -2077

# Row 113
$ws.Range("H113").Value = 500
$ws.Range("I113").Value = 500
$ws.Range("K113").Value = 500
$ws.Range("M113").Value = 1670

# Row 122
$ws.Range("H122").Value = 910.375
$ws.Range("J122").Value = 899.3333
$ws.Range("L122").Value = 2697.9999
$ws.Range("N122").Value = -7597.9999

$ws = $wb.Worksheets.Item("LTW")
# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# Row 46
$ws.Range("H46").Value = 3666.1667

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()

# Row 113
$ws.Range("H113").Value = 800
$ws.Range("J113").Value = 800
$ws.Range("L113").Value = 2400
$ws.Range("N113").Value = -6740

# Row 131
$ws.Range("H131").Value = 84900
$ws.Range("J131").Value = 84900
$ws.Range("L131").Value = 84900
$ws.Range("N131").Value = -94980

# Row 132
$ws.Range("H132").Value = 1999.8182
$ws.Range("I132").Value = 1300.2609
$ws.Range("J132").Value = 3608.8
$ws.Range("K132").Value = 3900.7827
$ws.Range("L132").Value = 10826.4
$ws.Range("M132").Value = -1370.7827
$ws.Range("N132").Value = -15886.4
